$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Existing transaction 8 (rows 11-13) gets its timestamp re-written with a
# tiny (sub-millisecond) precision correction.
$fixedDate = 45233.68449534722
$ws.Cells.Item(11, 2).Value = $fixedDate
$ws.Cells.Item(12, 2).Value = $fixedDate
$ws.Cells.Item(13, 2).Value = $fixedDate

# New transaction (Transaction_ID 9) split across 5 line items, same timestamp.
$newDate = 45234.50422506529

$rows = @(
    @(9, $newDate, 13,  13,  1, 7, "product 7"),
    @(9, $newDate, 15,  15,  1, 8, "product 8"),
    @(9, $newDate, 255, 255, 1, 9, "product 9"),
    @(9, $newDate, 7,   7,   1, 3, "product 3"),
    @(9, $newDate, 8,   4,   2, 2, "product 2")
)

$startRow = 14
for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $row = $rows[$i]

    $ws.Cells.Item($r, 1).Value = $row[0]

    $dateCell = $ws.Cells.Item($r, 2)
    $dateCell.Value = $row[1]
    $dateCell.NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]
}
